$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 09:20"

# Update country figures (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) for the countries whose numbers changed in
# this refresh. Rows below refer to the current (pre-sort) layout.

# Australia (row 21)
$ws.Cells.Item(21, 4).Value = 226
$ws.Cells.Item(21, 5).Value = 3727

# Pakistan (row 34)
$ws.Cells.Item(34, 2).Value = 1526
$ws.Cells.Item(34, 3).Value = 31
$ws.Cells.Item(34, 5).Value = 1484
$ws.Cells.Item(34, 6).Value = 11
$ws.Cells.Item(34, 7).Value = 1
$ws.Cells.Item(34, 8).Value = 13

# Tailandia (row 36)
$ws.Cells.Item(36, 5).Value = 1284
$ws.Cells.Item(36, 7).Value = 1
$ws.Cells.Item(36, 8).Value = 7

# Finlandia (row 40)
$ws.Cells.Item(40, 2).Value = 1218
$ws.Cells.Item(40, 3).Value = 51
$ws.Cells.Item(40, 4).Value = 10
$ws.Cells.Item(40, 5).Value = 1199
$ws.Cells.Item(40, 6).Value = 32
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 9

# Emiratos Arabes Unidos (row 64)
$ws.Cells.Item(64, 4).Value = 55
$ws.Cells.Item(64, 5).Value = 411

# Armenia (row 69)
$ws.Cells.Item(69, 2).Value = 424
$ws.Cells.Item(69, 3).Value = 17
$ws.Cells.Item(69, 4).Value = 30
$ws.Cells.Item(69, 5).Value = 391
$ws.Cells.Item(69, 6).Value = 6
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 3

# Letonia (row 74)
$ws.Cells.Item(74, 2).Value = 347
$ws.Cells.Item(74, 3).Value = 42
$ws.Cells.Item(74, 4).Value = 1
$ws.Cells.Item(74, 5).Value = 346
$ws.Cells.Item(74, 6).Value = 3
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0

# Bolivia (row 117)
$ws.Cells.Item(117, 6).Value = 3

# Re-sort the country table (rows 4-205) by "Casos totales" (column B) descending,
# which re-ranks the countries whose totals changed above.
$dataRange = $ws.Range("A4:H205")
$sortKey = $ws.Range("B4:B205")
$dataRange.Sort($sortKey, 2)
